# Apply the NATMI TPM recomputation update described in the commit.
# The MuSCs target-cluster rows are removed (receptor no longer detected there),
# which shrinks the table from 16 data rows (4x4 sender/target grid) to 12 rows
# (4 senders x 3 remaining targets: ECs, FAPs, Resolving-Mac), and every remaining
# metric is recomputed against the new (smaller) set of expressing clusters.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the four rows (14-17) that described the MuSCs target cluster;
# row count goes from 17 to 13 and Excel auto-shrinks the sheet dimension
# and the shared-string usage count accordingly.
$ws.Rows("14:17").Delete()

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 1.468507333333333
$ws.Range("H2").Value = 4.405521999999999
$ws.Range("I2").Value = 0.005118279455112885
$ws.Range("J2").Value = 0.005118279455112885
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.101223
$ws.Range("N2").Value = 0.303669
$ws.Range("O2").Value = 0.01313978392822635
$ws.Range("P2").Value = 0.01313978392822635
$ws.Range("Q2").Value = 0.148646717802
$ws.Range("R2").Value = 1.337820460218
$ws.Range("S2").Value = 0.00006725308612446341
$ws.Range("T2").Value = 0.00006725308612446341

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 1.468507333333333
$ws.Range("H3").Value = 4.405521999999999
$ws.Range("I3").Value = 0.005118279455112885
$ws.Range("J3").Value = 0.005118279455112885
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3740683333333334
$ws.Range("N3").Value = 1.122205
$ws.Range("O3").Value = 0.0485579075347673
$ws.Range("P3").Value = 0.0485579075347673
$ws.Range("Q3").Value = 0.5493220906677778
$ws.Range("R3").Value = 4.94389881601
$ws.Range("S3").Value = 0.0002485329405184707
$ws.Range("T3").Value = 0.0002485329405184707

# Row 4: ECs -> Resolving-Mac
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 1.468507333333333
$ws.Range("H4").Value = 4.405521999999999
$ws.Range("I4").Value = 0.005118279455112885
$ws.Range("J4").Value = 0.005118279455112885
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.22826
$ws.Range("N4").Value = 21.68478
$ws.Range("O4").Value = 0.9383023085370062
$ws.Range("P4").Value = 0.9383023085370064
$ws.Range("Q4").Value = 10.61475281724
$ws.Range("R4").Value = 95.53277535515998
$ws.Range("S4").Value = 0.00480249342846995
$ws.Range("T4").Value = 0.004802493428469951

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "ECs"
$ws.Range("G5").Value = 259.5505726666667
$ws.Range("H5").Value = 778.6517180000001
$ws.Range("I5").Value = 0.9046276674881553
$ws.Range("J5").Value = 0.9046276674881553
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.101223
$ws.Range("N5").Value = 0.303669
$ws.Range("O5").Value = 0.01313978392822635
$ws.Range("P5").Value = 0.01313978392822635
$ws.Range("Q5").Value = 26.27248761703801
$ws.Range("R5").Value = 236.452388553342
$ws.Range("S5").Value = 0.01188661208628975
$ws.Range("T5").Value = 0.01188661208628975

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 259.5505726666667
$ws.Range("H6").Value = 778.6517180000001
$ws.Range("I6").Value = 0.9046276674881553
$ws.Range("J6").Value = 0.9046276674881553
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3740683333333334
$ws.Range("N6").Value = 1.122205
$ws.Range("O6").Value = 0.0485579075347673
$ws.Range("P6").Value = 0.0485579075347673
$ws.Range("Q6").Value = 97.08965013313225
$ws.Range("R6").Value = 873.8068511981902
$ws.Range("S6").Value = 0.04392682663128206
$ws.Range("T6").Value = 0.04392682663128206

# Row 7: FAPs -> Resolving-Mac
$ws.Range("A7").Value = "FAPs"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("G7").Value = 259.5505726666667
$ws.Range("H7").Value = 778.6517180000001
$ws.Range("I7").Value = 0.9046276674881553
$ws.Range("J7").Value = 0.9046276674881553
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.22826
$ws.Range("N7").Value = 21.68478
$ws.Range("O7").Value = 0.9383023085370062
$ws.Range("P7").Value = 0.9383023085370064
$ws.Range("Q7").Value = 1876.09902238356
$ws.Range("R7").Value = 16884.89120145204
$ws.Range("S7").Value = 0.8488142287705833
$ws.Range("T7").Value = 0.8488142287705834

# Row 8: MuSCs -> ECs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("D8").Value = "ECs"
$ws.Range("G8").Value = 0.5890733333333333
$ws.Range("H8").Value = 1.76722
$ws.Range("I8").Value = 0.002053133730501083
$ws.Range("J8").Value = 0.002053133730501083
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.101223
$ws.Range("N8").Value = 0.303669
$ws.Range("O8").Value = 0.01313978392822635
$ws.Range("P8").Value = 0.01313978392822635
$ws.Range("Q8").Value = 0.05962777002000001
$ws.Range("R8").Value = 0.5366499301800001
$ws.Range("S8").Value = 0.00002697773359453754
$ws.Range("T8").Value = 0.00002697773359453754

# Row 9: MuSCs -> FAPs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("G9").Value = 0.5890733333333333
$ws.Range("H9").Value = 1.76722
$ws.Range("I9").Value = 0.002053133730501083
$ws.Range("J9").Value = 0.002053133730501083
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.3740683333333334
$ws.Range("N9").Value = 1.122205
$ws.Range("O9").Value = 0.0485579075347673
$ws.Range("P9").Value = 0.0485579075347673
$ws.Range("Q9").Value = 0.2203536800111111
$ws.Range("R9").Value = 1.9831831201
$ws.Range("S9").Value = 0.00009969587784218345
$ws.Range("T9").Value = 0.00009969587784218345

# Row 10: MuSCs -> Resolving-Mac
$ws.Range("A10").Value = "MuSCs"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("G10").Value = 0.5890733333333333
$ws.Range("H10").Value = 1.76722
$ws.Range("I10").Value = 0.002053133730501083
$ws.Range("J10").Value = 0.002053133730501083
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.22826
$ws.Range("N10").Value = 21.68478
$ws.Range("O10").Value = 0.9383023085370062
$ws.Range("P10").Value = 0.9383023085370064
$ws.Range("Q10").Value = 4.2579752124
$ws.Range("R10").Value = 38.3217769116
$ws.Range("S10").Value = 0.001926460119064362
$ws.Range("T10").Value = 0.001926460119064362

# Row 11: Resolving-Mac -> ECs
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("D11").Value = "ECs"
$ws.Range("G11").Value = 25.306101
$ws.Range("H11").Value = 75.918303
$ws.Range("I11").Value = 0.0882009193262308
$ws.Range("J11").Value = 0.0882009193262308
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.101223
$ws.Range("N11").Value = 0.303669
$ws.Range("O11").Value = 0.01313978392822635
$ws.Range("P11").Value = 0.01313978392822635
$ws.Range("Q11").Value = 2.561559461523
$ws.Range("R11").Value = 23.054035153707
$ws.Range("S11").Value = 0.001158941022217596
$ws.Range("T11").Value = 0.001158941022217596

# Row 12: Resolving-Mac -> FAPs
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("D12").Value = "FAPs"
$ws.Range("G12").Value = 25.306101
$ws.Range("H12").Value = 75.918303
$ws.Range("I12").Value = 0.0882009193262308
$ws.Range("J12").Value = 0.0882009193262308
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.3740683333333334
$ws.Range("N12").Value = 1.122205
$ws.Range("O12").Value = 0.0485579075347673
$ws.Range("P12").Value = 0.0485579075347673
$ws.Range("Q12").Value = 9.466211024235001
$ws.Range("R12").Value = 85.195899218115
$ws.Range("S12").Value = 0.004282852085124586
$ws.Range("T12").Value = 0.004282852085124586

# Row 13: Resolving-Mac -> Resolving-Mac
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("G13").Value = 25.306101
$ws.Range("H13").Value = 75.918303
$ws.Range("I13").Value = 0.0882009193262308
$ws.Range("J13").Value = 0.0882009193262308
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 7.22826
$ws.Range("N13").Value = 21.68478
$ws.Range("O13").Value = 0.9383023085370062
$ws.Range("P13").Value = 0.9383023085370064
$ws.Range("Q13").Value = 182.91907761426
$ws.Range("R13").Value = 1646.27169852834
$ws.Range("S13").Value = 0.0827591262188886
$ws.Range("T13").Value = 0.08275912621888862
